# Leave Card update: add two SL(1-0-0) single-day entries (rows 90 & 91),
# insert two new rows for an SL(3-0-0) multi-day entry (new rows 92 & 93),
# and update the saved absolute path + view selection to match the author's
# re-save of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$tbl = $ws.ListObjects.Item("Table1")

# --- Fill in the two existing blank rows (90 & 91) with SL(1-0-0) entries ---
$ws.Range("B90").Value = "SL(1-0-0)"
$ws.Range("C90").Value2 = 1.25
$ws.Range("H90").Value2 = 1
$ws.Range("K90").NumberFormat = "m/d/yyyy"
$ws.Range("K90").Value2 = 45142

$ws.Range("B91").Value = "SL(1-0-0)"
$ws.Range("C91").Value2 = 1.25
$ws.Range("H91").Value2 = 1
$ws.Range("K91").NumberFormat = "m/d/yyyy"
$ws.Range("K91").Value2 = 45163

# --- Insert two new table rows right after row 91 (i.e. before the old row 92) ---
# ListRows index is 1-based relative to the table body (row 9 == index 1),
# so sheet row 92 is table-row index 84.
$tbl.ListRows.Add(84) | Out-Null
$tbl.ListRows.Add(84) | Out-Null

# The PERIOD column isn't a calculated table column, but the insert still
# carries a fill-down date formula into the two fresh rows; the real edit
# left those two PERIOD cells blank, so clear them back out.
$ws.Range("A92").ClearContents()
$ws.Range("A93").ClearContents()

# Row 94 (the old row 92) keeps its PERIOD formula, but since the two rows
# above it are now blank it chains straight back to A91 instead of A93.
$ws.Range("A94").Formula = "=EDATE(A91,1)"

# New row 92: SL(3-0-0), spanning 8/30 - 9/1/2023 (3 hours)
$ws.Range("B92").Value = "SL(3-0-0)"
$ws.Range("H92").Value2 = 3
$ws.Range("K92").NumberFormat = "m/d/yyyy"
$ws.Range("K92").Value = "8/30 - 9/1/2023"

# New row 93: SL(1-0-0) (1 hour)
$ws.Range("B93").Value = "SL(1-0-0)"
$ws.Range("H93").Value2 = 1
$ws.Range("K93").NumberFormat = "m/d/yyyy"
$ws.Range("K93").Value2 = 45190

# --- Match the saved absolute path recorded by Excel on re-save ---
$wb.Windows.Item(1).Caption = $wb.Windows.Item(1).Caption

# --- Restore/refresh the view state (scroll position + active selection) ---
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 81
$ws.Range("K93").Select()

$wb.RecalcBeforeSave = $true
